# Actualización automatica mar abr  6 17:33:39 CEST 2021
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (row 3: sdmx dim tags; row 4: dim/medida; row 5: xsd types) ---
$ws.Range("B3").Value = "iaest-measure:orden"
$ws.Range("G3").Value = "iaest-measure:siglas-agrupada"
$ws.Range("J3").Value = "iaest-measure:siglas"

$ws.Range("B4").Value = "medida"
$ws.Range("G4").Value = "medida"
$ws.Range("J4").Value = "medida"

$ws.Range("B5").Value = "xsd:int"
$ws.Range("G5").Value = "xsd:string"
$ws.Range("J5").Value = "xsd:string"

# --- Normalize formatting: cells that carried an explicit Arial-only style
# collapse onto the default (Arial 10 / implicit) style used by the rest
# of the sheet. ---
$ws.Range("J1").Font.Name = "Arial"
$ws.Range("J2").Font.Name = "Arial"
$ws.Range("B3").Font.Name = "Arial"
$ws.Range("J3").Font.Name = "Arial"
$ws.Range("B4").Font.Name = "Arial"
$ws.Range("G4").Font.Name = "Arial"
$ws.Range("J4").Font.Name = "Arial"
$ws.Range("B5").Font.Name = "Arial"
$ws.Range("G5").Font.Name = "Arial"
$ws.Range("J5").Font.Name = "Arial"

# --- Drop the trailing mapping-file row; it is no longer part of the table ---
$ws.Rows(6).Delete()
